$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "completed tank titrations 0224" -- append the 2022-02-24 titration
# result as a new row (row 60) below the existing data table.
$ws.Range("A60").Value = 20220224
$ws.Range("B60").Value = 2222.0050000000001
$ws.Range("C60").Value = 2224.4699999999998
$ws.Range("D60").Formula = "=100*(B60-C60)/C60"
$ws.Range("E60").Value = 180
$ws.Range("F60").Value = "CRM OPENED 20220118"

# Match the author's window position after entering the new row: scrolled
# down so row 45 is at the top, with the cell just past the new row active.
$win = $excel.ActiveWindow
$win.ScrollRow = 45
$win.ScrollColumn = 1
$ws.Range("I60").Select()
